$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,7).Value2 = 4.564056666666667
$ws.Cells.Item(2,8).Value2 = 13.69217
$ws.Cells.Item(2,9).Value2 = 0.3012303010600196
$ws.Cells.Item(2,10).Value2 = 0.3202575519860646
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,13).Value2 = 3.946674666666667
$ws.Cells.Item(2,14).Value2 = 11.840024
$ws.Cells.Item(2,15).Value2 = 0.008082287850277156
$ws.Cells.Item(2,16).Value2 = 0.008317683527585098
$ws.Cells.Item(2,17).Value2 = 18.01284682356444
$ws.Cells.Item(2,18).Value2 = 162.11562141208
$ws.Cells.Item(2,19).Value2 = 0.002434630002392727
$ws.Cells.Item(2,20).Value2 = 0.002663800964739218
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,7).Value2 = 4.564056666666667
$ws.Cells.Item(3,8).Value2 = 13.69217
$ws.Cells.Item(3,9).Value2 = 0.3012303010600196
$ws.Cells.Item(3,10).Value2 = 0.3202575519860646
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,13).Value2 = 185.8027443333333
$ws.Cells.Item(3,14).Value2 = 557.408233
$ws.Cells.Item(3,15).Value2 = 0.3805003933455167
$ws.Cells.Item(3,16).Value2 = 0.3915824222792467
$ws.Cells.Item(3,17).Value2 = 848.0142539595122
$ws.Cells.Item(3,18).Value2 = 7632.12828563561
$ws.Cells.Item(3,19).Value2 = 0.1146182480409259
$ws.Cells.Item(3,20).Value2 = 0.125407227959925
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,7).Value2 = 4.564056666666667
$ws.Cells.Item(4,8).Value2 = 13.69217
$ws.Cells.Item(4,9).Value2 = 0.3012303010600196
$ws.Cells.Item(4,10).Value2 = 0.3202575519860646
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,13).Value2 = 117.3394243333333
$ws.Cells.Item(4,14).Value2 = 352.018273
$ws.Cells.Item(4,15).Value2 = 0.2402962199184265
$ws.Cells.Item(4,16).Value2 = 0.2472948188906589
$ws.Cells.Item(4,17).Value2 = 535.5437818913789
$ws.Cells.Item(4,18).Value2 = 4819.894037022411
$ws.Cells.Item(4,19).Value2 = 0.07238450266961229
$ws.Cells.Item(4,20).Value2 = 0.0791980333167596
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,7).Value2 = 4.564056666666667
$ws.Cells.Item(5,8).Value2 = 13.69217
$ws.Cells.Item(5,9).Value2 = 0.3012303010600196
$ws.Cells.Item(5,10).Value2 = 0.3202575519860646
$ws.Cells.Item(5,11).Value2 = 3
$ws.Cells.Item(5,13).Value2 = 139.7641523333333
$ws.Cells.Item(5,14).Value2 = 419.292457
$ws.Cells.Item(5,15).Value2 = 0.2862192112890951
$ws.Cells.Item(5,16).Value2 = 0.2945553119511906
$ws.Cells.Item(5,17).Value2 = 637.8915112179656
$ws.Cells.Item(5,18).Value2 = 5741.02360096169
$ws.Cells.Item(5,19).Value2 = 0.08621789918577548
$ws.Cells.Item(5,20).Value2 = 0.09433356312997988
$ws.Cells.Item(6,5).Value2 = 3
$ws.Cells.Item(6,7).Value2 = 4.564056666666667
$ws.Cells.Item(6,8).Value2 = 13.69217
$ws.Cells.Item(6,9).Value2 = 0.3012303010600196
$ws.Cells.Item(6,10).Value2 = 0.3202575519860646
$ws.Cells.Item(6,11).Value2 = 2
$ws.Cells.Item(6,13).Value2 = 41.458574
$ws.Cells.Item(6,14).Value2 = 82.917148
$ws.Cells.Item(6,15).Value2 = 0.0849018875966847
$ws.Cells.Item(6,16).Value2 = 0.05824976335131885
$ws.Cells.Item(6,17).Value2 = 189.2192810551933
$ws.Cells.Item(6,18).Value2 = 1135.31568633116
$ws.Cells.Item(6,19).Value2 = 0.02557502116131328
$ws.Cells.Item(6,20).Value2 = 0.01865492661466096
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,7).Value2 = 2.66687
$ws.Cells.Item(7,8).Value2 = 8.00061
$ws.Cells.Item(7,9).Value2 = 0.1760149164788199
$ws.Cells.Item(7,10).Value2 = 0.187132921443075
$ws.Cells.Item(7,11).Value2 = 3
$ws.Cells.Item(7,13).Value2 = 3.946674666666667
$ws.Cells.Item(7,14).Value2 = 11.840024
$ws.Cells.Item(7,15).Value2 = 0.008082287850277156
$ws.Cells.Item(7,16).Value2 = 0.008317683527585098
$ws.Cells.Item(7,17).Value2 = 10.52526826829333
$ws.Cells.Item(7,18).Value2 = 94.72741441464
$ws.Cells.Item(7,19).Value2 = 0.001422603220924315
$ws.Cells.Item(7,20).Value2 = 0.001556512418155941
$ws.Cells.Item(8,5).Value2 = 3
$ws.Cells.Item(8,7).Value2 = 2.66687
$ws.Cells.Item(8,8).Value2 = 8.00061
$ws.Cells.Item(8,9).Value2 = 0.1760149164788199
$ws.Cells.Item(8,10).Value2 = 0.187132921443075
$ws.Cells.Item(8,11).Value2 = 3
$ws.Cells.Item(8,13).Value2 = 185.8027443333333
$ws.Cells.Item(8,14).Value2 = 557.408233
$ws.Cells.Item(8,15).Value2 = 0.3805003933455167
$ws.Cells.Item(8,16).Value2 = 0.3915824222792467
$ws.Cells.Item(8,17).Value2 = 495.5117647802366
$ws.Cells.Item(8,18).Value2 = 4459.60588302213
$ws.Cells.Item(8,19).Value2 = 0.06697374495486923
$ws.Cells.Item(8,20).Value2 = 0.07327796266687128
$ws.Cells.Item(9,5).Value2 = 3
$ws.Cells.Item(9,7).Value2 = 2.66687
$ws.Cells.Item(9,8).Value2 = 8.00061
$ws.Cells.Item(9,9).Value2 = 0.1760149164788199
$ws.Cells.Item(9,10).Value2 = 0.187132921443075
$ws.Cells.Item(9,11).Value2 = 3
$ws.Cells.Item(9,13).Value2 = 117.3394243333333
$ws.Cells.Item(9,14).Value2 = 352.018273
$ws.Cells.Item(9,15).Value2 = 0.2402962199184265
$ws.Cells.Item(9,16).Value2 = 0.2472948188906589
$ws.Cells.Item(9,17).Value2 = 312.9289905718367
$ws.Cells.Item(9,18).Value2 = 2816.36091514653
$ws.Cells.Item(9,19).Value2 = 0.04229571907911797
$ws.Cells.Item(9,20).Value2 = 0.04627700191674512
$ws.Cells.Item(10,5).Value2 = 3
$ws.Cells.Item(10,7).Value2 = 2.66687
$ws.Cells.Item(10,8).Value2 = 8.00061
$ws.Cells.Item(10,9).Value2 = 0.1760149164788199
$ws.Cells.Item(10,10).Value2 = 0.187132921443075
$ws.Cells.Item(10,11).Value2 = 3
$ws.Cells.Item(10,13).Value2 = 139.7641523333333
$ws.Cells.Item(10,14).Value2 = 419.292457
$ws.Cells.Item(10,15).Value2 = 0.2862192112890951
$ws.Cells.Item(10,16).Value2 = 0.2945553119511906
$ws.Cells.Item(10,17).Value2 = 372.7328249331966
$ws.Cells.Item(10,18).Value2 = 3354.59542439877
$ws.Cells.Item(10,19).Value2 = 0.05037885056968377
$ws.Cells.Item(10,20).Value2 = 0.05512099605200259
$ws.Cells.Item(11,5).Value2 = 3
$ws.Cells.Item(11,7).Value2 = 2.66687
$ws.Cells.Item(11,8).Value2 = 8.00061
$ws.Cells.Item(11,9).Value2 = 0.1760149164788199
$ws.Cells.Item(11,10).Value2 = 0.187132921443075
$ws.Cells.Item(11,11).Value2 = 2
$ws.Cells.Item(11,13).Value2 = 41.458574
$ws.Cells.Item(11,14).Value2 = 82.917148
$ws.Cells.Item(11,15).Value2 = 0.0849018875966847
$ws.Cells.Item(11,16).Value2 = 0.05824976335131885
$ws.Cells.Item(11,17).Value2 = 110.56462724338
$ws.Cells.Item(11,18).Value2 = 663.38776346028
$ws.Cells.Item(11,19).Value2 = 0.01494399865422461
$ws.Cells.Item(11,20).Value2 = 0.01090044838930006
$ws.Cells.Item(12,5).Value2 = 3
$ws.Cells.Item(12,7).Value2 = 2.718648333333333
$ws.Cells.Item(12,8).Value2 = 8.155945
$ws.Cells.Item(12,9).Value2 = 0.1794323155335466
$ws.Cells.Item(12,10).Value2 = 0.1907661809510825
$ws.Cells.Item(12,11).Value2 = 3
$ws.Cells.Item(12,13).Value2 = 3.946674666666667
$ws.Cells.Item(12,14).Value2 = 11.840024
$ws.Cells.Item(12,15).Value2 = 0.008082287850277156
$ws.Cells.Item(12,16).Value2 = 0.008317683527585098
$ws.Cells.Item(12,17).Value2 = 10.72962050474222
$ws.Cells.Item(12,18).Value2 = 96.56658454267999
$ws.Cells.Item(12,19).Value2 = 0.001450223623783881
$ws.Cells.Item(12,20).Value2 = 0.001586732720917137
$ws.Cells.Item(13,5).Value2 = 3
$ws.Cells.Item(13,7).Value2 = 2.718648333333333
$ws.Cells.Item(13,8).Value2 = 8.155945
$ws.Cells.Item(13,9).Value2 = 0.1794323155335466
$ws.Cells.Item(13,10).Value2 = 0.1907661809510825
$ws.Cells.Item(13,11).Value2 = 3
$ws.Cells.Item(13,13).Value2 = 185.8027443333333
$ws.Cells.Item(13,14).Value2 = 557.408233
$ws.Cells.Item(13,15).Value2 = 0.3805003933455167
$ws.Cells.Item(13,16).Value2 = 0.3915824222792467
$ws.Cells.Item(13,17).Value2 = 505.132321210576
$ws.Cells.Item(13,18).Value2 = 4546.190890895185
$ws.Cells.Item(13,19).Value2 = 0.06827406663941135
$ws.Cells.Item(13,20).Value2 = 0.07470068322578596
$ws.Cells.Item(14,5).Value2 = 3
$ws.Cells.Item(14,7).Value2 = 2.718648333333333
$ws.Cells.Item(14,8).Value2 = 8.155945
$ws.Cells.Item(14,9).Value2 = 0.1794323155335466
$ws.Cells.Item(14,10).Value2 = 0.1907661809510825
$ws.Cells.Item(14,11).Value2 = 3
$ws.Cells.Item(14,13).Value2 = 117.3394243333333
$ws.Cells.Item(14,14).Value2 = 352.018273
$ws.Cells.Item(14,15).Value2 = 0.2402962199184265
$ws.Cells.Item(14,16).Value2 = 0.2472948188906589
$ws.Cells.Item(14,17).Value2 = 319.0046303981094
$ws.Cells.Item(14,18).Value2 = 2871.041673582985
$ws.Cells.Item(14,19).Value2 = 0.04311690715392161
$ws.Cells.Item(14,20).Value2 = 0.0471754881687606
$ws.Cells.Item(15,5).Value2 = 3
$ws.Cells.Item(15,7).Value2 = 2.718648333333333
$ws.Cells.Item(15,8).Value2 = 8.155945
$ws.Cells.Item(15,9).Value2 = 0.1794323155335466
$ws.Cells.Item(15,10).Value2 = 0.1907661809510825
$ws.Cells.Item(15,11).Value2 = 3
$ws.Cells.Item(15,13).Value2 = 139.7641523333333
$ws.Cells.Item(15,14).Value2 = 419.292457
$ws.Cells.Item(15,15).Value2 = 0.2862192112890951
$ws.Cells.Item(15,16).Value2 = 0.2945553119511906
$ws.Cells.Item(15,17).Value2 = 379.9695798007627
$ws.Cells.Item(15,18).Value2 = 3419.726218206865
$ws.Cells.Item(15,19).Value2 = 0.05135697583178776
$ws.Cells.Item(15,20).Value2 = 0.05619119193978336
$ws.Cells.Item(16,5).Value2 = 3
$ws.Cells.Item(16,7).Value2 = 2.718648333333333
$ws.Cells.Item(16,8).Value2 = 8.155945
$ws.Cells.Item(16,9).Value2 = 0.1794323155335466
$ws.Cells.Item(16,10).Value2 = 0.1907661809510825
$ws.Cells.Item(16,11).Value2 = 2
$ws.Cells.Item(16,13).Value2 = 41.458574
$ws.Cells.Item(16,14).Value2 = 82.917148
$ws.Cells.Item(16,15).Value2 = 0.0849018875966847
$ws.Cells.Item(16,16).Value2 = 0.05824976335131885
$ws.Cells.Item(16,17).Value2 = 112.7112831074766
$ws.Cells.Item(16,18).Value2 = 676.2676986448599
$ws.Cells.Item(16,19).Value2 = 0.01523414228464204
$ws.Cells.Item(16,20).Value2 = 0.01111208489583542
$ws.Cells.Item(17,5).Value2 = 3
$ws.Cells.Item(17,7).Value2 = 2.501273333333333
$ws.Cells.Item(17,8).Value2 = 7.50382
$ws.Cells.Item(17,9).Value2 = 0.1650854435564411
$ws.Cells.Item(17,10).Value2 = 0.1755130869499919
$ws.Cells.Item(17,11).Value2 = 3
$ws.Cells.Item(17,13).Value2 = 3.946674666666667
$ws.Cells.Item(17,14).Value2 = 11.840024
$ws.Cells.Item(17,15).Value2 = 0.008082287850277156
$ws.Cells.Item(17,16).Value2 = 0.008317683527585098
$ws.Cells.Item(17,17).Value2 = 9.871712099075555
$ws.Cells.Item(17,18).Value2 = 88.84540889168
$ws.Cells.Item(17,19).Value2 = 0.001334268074713839
$ws.Cells.Item(17,20).Value2 = 0.001459862312199559
$ws.Cells.Item(18,5).Value2 = 3
$ws.Cells.Item(18,7).Value2 = 2.501273333333333
$ws.Cells.Item(18,8).Value2 = 7.50382
$ws.Cells.Item(18,9).Value2 = 0.1650854435564411
$ws.Cells.Item(18,10).Value2 = 0.1755130869499919
$ws.Cells.Item(18,11).Value2 = 3
$ws.Cells.Item(18,13).Value2 = 185.8027443333333
$ws.Cells.Item(18,14).Value2 = 557.408233
$ws.Cells.Item(18,15).Value2 = 0.3805003933455167
$ws.Cells.Item(18,16).Value2 = 0.3915824222792467
$ws.Cells.Item(18,17).Value2 = 464.7434496611177
$ws.Cells.Item(18,18).Value2 = 4182.69104695006
$ws.Cells.Item(18,19).Value2 = 0.06281507620884494
$ws.Cells.Item(18,20).Value2 = 0.06872783972958589
$ws.Cells.Item(19,5).Value2 = 3
$ws.Cells.Item(19,7).Value2 = 2.501273333333333
$ws.Cells.Item(19,8).Value2 = 7.50382
$ws.Cells.Item(19,9).Value2 = 0.1650854435564411
$ws.Cells.Item(19,10).Value2 = 0.1755130869499919
$ws.Cells.Item(19,11).Value2 = 3
$ws.Cells.Item(19,13).Value2 = 117.3394243333333
$ws.Cells.Item(19,14).Value2 = 352.018273
$ws.Cells.Item(19,15).Value2 = 0.2402962199184265
$ws.Cells.Item(19,16).Value2 = 0.2472948188906589
$ws.Cells.Item(19,17).Value2 = 293.4979730336511
$ws.Cells.Item(19,18).Value2 = 2641.48175730286
$ws.Cells.Item(19,19).Value2 = 0.03966940805016955
$ws.Cells.Item(19,20).Value2 = 0.04340347705023871
$ws.Cells.Item(20,5).Value2 = 3
$ws.Cells.Item(20,7).Value2 = 2.501273333333333
$ws.Cells.Item(20,8).Value2 = 7.50382
$ws.Cells.Item(20,9).Value2 = 0.1650854435564411
$ws.Cells.Item(20,10).Value2 = 0.1755130869499919
$ws.Cells.Item(20,11).Value2 = 3
$ws.Cells.Item(20,13).Value2 = 139.7641523333333
$ws.Cells.Item(20,14).Value2 = 419.292457
$ws.Cells.Item(20,15).Value2 = 0.2862192112890951
$ws.Cells.Item(20,16).Value2 = 0.2945553119511906
$ws.Cells.Item(20,17).Value2 = 349.5883471873044
$ws.Cells.Item(20,18).Value2 = 3146.29512468574
$ws.Cells.Item(20,19).Value2 = 0.04725062545003499
$ws.Cells.Item(20,20).Value2 = 0.0516983120780713
$ws.Cells.Item(21,5).Value2 = 3
$ws.Cells.Item(21,7).Value2 = 2.501273333333333
$ws.Cells.Item(21,8).Value2 = 7.50382
$ws.Cells.Item(21,9).Value2 = 0.1650854435564411
$ws.Cells.Item(21,10).Value2 = 0.1755130869499919
$ws.Cells.Item(21,11).Value2 = 2
$ws.Cells.Item(21,13).Value2 = 41.458574
$ws.Cells.Item(21,14).Value2 = 82.917148
$ws.Cells.Item(21,15).Value2 = 0.0849018875966847
$ws.Cells.Item(21,16).Value2 = 0.05824976335131885
$ws.Cells.Item(21,17).Value2 = 103.6992255842267
$ws.Cells.Item(21,18).Value2 = 622.19535350536
$ws.Cells.Item(21,19).Value2 = 0.0140160657726778
$ws.Cells.Item(21,20).Value2 = 0.01022359577989648
$ws.Cells.Item(22,5).Value2 = 2
$ws.Cells.Item(22,7).Value2 = 2.700538
$ws.Cells.Item(22,8).Value2 = 5.401076
$ws.Cells.Item(22,9).Value2 = 0.1782370233711727
$ws.Cells.Item(22,10).Value2 = 0.1263302586697861
$ws.Cells.Item(22,11).Value2 = 3
$ws.Cells.Item(22,13).Value2 = 3.946674666666667
$ws.Cells.Item(22,14).Value2 = 11.840024
$ws.Cells.Item(22,15).Value2 = 0.008082287850277156
$ws.Cells.Item(22,16).Value2 = 0.008317683527585098
$ws.Cells.Item(22,17).Value2 = 10.65814491097067
$ws.Cells.Item(22,18).Value2 = 63.948869465824
$ws.Cells.Item(22,19).Value2 = 0.001440562928462395
$ws.Cells.Item(22,20).Value2 = 0.001050775111573244
$ws.Cells.Item(23,5).Value2 = 2
$ws.Cells.Item(23,7).Value2 = 2.700538
$ws.Cells.Item(23,8).Value2 = 5.401076
$ws.Cells.Item(23,9).Value2 = 0.1782370233711727
$ws.Cells.Item(23,10).Value2 = 0.1263302586697861
$ws.Cells.Item(23,11).Value2 = 3
$ws.Cells.Item(23,13).Value2 = 185.8027443333333
$ws.Cells.Item(23,14).Value2 = 557.408233
$ws.Cells.Item(23,15).Value2 = 0.3805003933455167
$ws.Cells.Item(23,16).Value2 = 0.3915824222792467
$ws.Cells.Item(23,17).Value2 = 501.7673715764513
$ws.Cells.Item(23,18).Value2 = 3010.604229458708
$ws.Cells.Item(23,19).Value2 = 0.06781925750146525
$ws.Cells.Item(23,20).Value2 = 0.04946870869707866
$ws.Cells.Item(24,5).Value2 = 2
$ws.Cells.Item(24,7).Value2 = 2.700538
$ws.Cells.Item(24,8).Value2 = 5.401076
$ws.Cells.Item(24,9).Value2 = 0.1782370233711727
$ws.Cells.Item(24,10).Value2 = 0.1263302586697861
$ws.Cells.Item(24,11).Value2 = 3
$ws.Cells.Item(24,13).Value2 = 117.3394243333333
$ws.Cells.Item(24,14).Value2 = 352.018273
$ws.Cells.Item(24,15).Value2 = 0.2402962199184265
$ws.Cells.Item(24,16).Value2 = 0.2472948188906589
$ws.Cells.Item(24,17).Value2 = 316.8795743102913
$ws.Cells.Item(24,18).Value2 = 1901.277445861748
$ws.Cells.Item(24,19).Value2 = 0.04282968296560503
$ws.Cells.Item(24,20).Value2 = 0.03124081843815484
$ws.Cells.Item(25,5).Value2 = 2
$ws.Cells.Item(25,7).Value2 = 2.700538
$ws.Cells.Item(25,8).Value2 = 5.401076
$ws.Cells.Item(25,9).Value2 = 0.1782370233711727
$ws.Cells.Item(25,10).Value2 = 0.1263302586697861
$ws.Cells.Item(25,11).Value2 = 3
$ws.Cells.Item(25,13).Value2 = 139.7641523333333
$ws.Cells.Item(25,14).Value2 = 419.292457
$ws.Cells.Item(25,15).Value2 = 0.2862192112890951
$ws.Cells.Item(25,16).Value2 = 0.2945553119511906
$ws.Cells.Item(25,17).Value2 = 377.4384044139553
$ws.Cells.Item(25,18).Value2 = 2264.630426483732
$ws.Cells.Item(25,19).Value2 = 0.05101486025181305
$ws.Cells.Item(25,20).Value2 = 0.03721124875135345
$ws.Cells.Item(26,5).Value2 = 2
$ws.Cells.Item(26,7).Value2 = 2.700538
$ws.Cells.Item(26,8).Value2 = 5.401076
$ws.Cells.Item(26,9).Value2 = 0.1782370233711727
$ws.Cells.Item(26,10).Value2 = 0.1263302586697861
$ws.Cells.Item(26,11).Value2 = 2
$ws.Cells.Item(26,13).Value2 = 41.458574
$ws.Cells.Item(26,14).Value2 = 82.917148
$ws.Cells.Item(26,15).Value2 = 0.0849018875966847
$ws.Cells.Item(26,16).Value2 = 0.05824976335131885
$ws.Cells.Item(26,17).Value2 = 111.960454512812
$ws.Cells.Item(26,18).Value2 = 447.841818051248
$ws.Cells.Item(26,19).Value2 = 0.01513265972382697
$ws.Cells.Item(26,20).Value2 = 0.007358707671625938
